$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 667291.75
$ws.Range("I38").Value = 1428760.2
$ws.Range("J38").Value = 1006.75
$ws.Range("K38").Value = 4286280.6
$ws.Range("L38").Value = 3020.25
$ws.Range("M38").Value = -4285908.6
$ws.Range("N38").Value = -3764.25

$ws.Range("H112").Value = 1665.8064
$ws.Range("J112").Value = 1704.6666
$ws.Range("L112").Value = 5113.9998
$ws.Range("N112").Value = -7329.9998

$ws.Range("H132").Value = 112803.586
$ws.Range("I132").Value = 1350.1803
$ws.Range("K132").Value = 4050.5409
$ws.Range("M132").Value = -1520.5409

$ws.Range("H137").Value = 21779.021
$ws.Range("I137").Value = 27881.027
$ws.Range("J137").Value = 1254.091
$ws.Range("K137").Value = 83643.08099999999
$ws.Range("L137").Value = 3762.273
$ws.Range("M137").Value = -81093.08099999999
$ws.Range("N137").Value = -8862.272999999999

$ws.Range("H138").Value = 1412.23
$ws.Range("I138").Value = 730.67444
$ws.Range("J138").Value = 1926.386
$ws.Range("K138").Value = 2192.02332
$ws.Range("L138").Value = 5779.157999999999
$ws.Range("M138").Value = 2947.97668
$ws.Range("N138").Value = -16059.158

$ws.Range("H141").Value = 1745.6938
$ws.Range("I141").Value = 997.15625
$ws.Range("J141").Value = 3154.7058
$ws.Range("K141").Value = 2991.46875
$ws.Range("L141").Value = 9464.117400000001
$ws.Range("M141").Value = 2188.53125
$ws.Range("N141").Value = -19824.1174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1320948.8
$ws.Range("I32").Value = 1567760
$ws.Range("J32").Value = 4621.5835
$ws.Range("K32").Value = 1567760
$ws.Range("L32").Value = 4621.5835
$ws.Range("M32").Value = -1567473
$ws.Range("N32").Value = -5195.5835

$ws.Range("H61").Value = 687.0984
$ws.Range("I61").Value = 691.96
$ws.Range("J61").Value = 665
$ws.Range("K61").Value = 691.96
$ws.Range("L61").Value = 665
$ws.Range("M61").Value = -479.96
$ws.Range("N61").Value = -1089

$ws.Range("H74").Value = 42092.918
$ws.Range("I74").Value = 70363.484
$ws.Range("J74").Value = 1100.6
$ws.Range("K74").Value = 70363.484
$ws.Range("L74").Value = 1100.6
$ws.Range("M74").Value = -69489.484
$ws.Range("N74").Value = -2848.6

$ws.Range("H77").Value = 42092.918
$ws.Range("I77").Value = 70363.484
$ws.Range("J77").Value = 1100.6
$ws.Range("K77").Value = 351817.42
$ws.Range("L77").Value = 5503
$ws.Range("M77").Value = -347449.42
$ws.Range("N77").Value = -14239

$ws.Range("H132").Value = 1871121.5
$ws.Range("I132").Value = 2431965.2
$ws.Range("J132").Value = 562485.9
$ws.Range("K132").Value = 7295895.600000001
$ws.Range("L132").Value = 1687457.7
$ws.Range("M132").Value = -7293365.600000001
$ws.Range("N132").Value = -1692517.7

$ws.Range("H136").Value = 687.0984
$ws.Range("I136").Value = 691.96
$ws.Range("J136").Value = 665
$ws.Range("K136").Value = 2075.88
$ws.Range("L136").Value = 1995
$ws.Range("M136").Value = 474.1199999999999
$ws.Range("N136").Value = -7095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25286.826
$ws.Range("I134").Value = 1087.9445
$ws.Range("J134").Value = 112402.8
$ws.Range("K134").Value = 3263.8335
$ws.Range("L134").Value = 337208.4
$ws.Range("M134").Value = -728.8335000000002
$ws.Range("N134").Value = -342278.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27879.076
$ws.Range("I31").Value = 34528.49
$ws.Range("K31").Value = 34528.49
$ws.Range("M31").Value = -34233.49

$ws.Range("H34").Value = 27879.076
$ws.Range("I34").Value = 34528.49
$ws.Range("K34").Value = 34528.49
$ws.Range("M34").Value = -34326.49

$ws.Range("H58").Value = 1832.0896
$ws.Range("I58").Value = 480.30612
$ws.Range("J58").Value = 5511.9443
$ws.Range("K58").Value = 480.30612
$ws.Range("L58").Value = 5511.9443
$ws.Range("M58").Value = -277.30612
$ws.Range("N58").Value = -5917.9443

$ws.Range("H132").Value = 1084.9143
$ws.Range("I132").Value = 841.5357
$ws.Range("J132").Value = 2058.4285
$ws.Range("K132").Value = 2524.6071
$ws.Range("L132").Value = 6175.2855
$ws.Range("M132").Value = 5.392899999999827
$ws.Range("N132").Value = -11235.2855

$ws.Range("H134").Value = 1280.381
$ws.Range("I134").Value = 1119.1562
$ws.Range("J134").Value = 1796.3
$ws.Range("K134").Value = 3357.4686
$ws.Range("L134").Value = 5388.9
$ws.Range("M134").Value = -822.4685999999997
$ws.Range("N134").Value = -10458.9

$ws.Range("H136").Value = 1832.0896
$ws.Range("I136").Value = 480.30612
$ws.Range("J136").Value = 5511.9443
$ws.Range("K136").Value = 1440.91836
$ws.Range("L136").Value = 16535.8329
$ws.Range("M136").Value = 1109.08164
$ws.Range("N136").Value = -21635.8329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14535725
$ws.Range("J131").Value = 16892814
$ws.Range("L131").Value = 50678442
$ws.Range("N131").Value = -50688522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 47337.332
$ws.Range("I20").Value = 2000
$ws.Range("K20").Value = 2000
$ws.Range("M20").Value = -1755

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 180868.08
$ws.Range("I132").Value = 44949.434
$ws.Range("J132").Value = 478594.62
$ws.Range("K132").Value = 134848.302
$ws.Range("L132").Value = 1435783.86
$ws.Range("M132").Value = -132318.302
$ws.Range("N132").Value = -1440843.86

$ws.Range("H136").Value = 190066.16
$ws.Range("I136").Value = 286808.78
$ws.Range("J136").Value = 1955.5
$ws.Range("K136").Value = 860426.3400000001
$ws.Range("L136").Value = 5866.5
$ws.Range("M136").Value = -857876.3400000001
$ws.Range("N136").Value = -10966.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 29404.4
$ws.Range("J23").Value = 29404.4
$ws.Range("L23").Value = 29404.4
$ws.Range("N23").Value = -29862.4

$ws.Range("H132").Value = 3127.6
$ws.Range("I132").Value = 562.75
$ws.Range("J132").Value = 13387
$ws.Range("K132").Value = 1688.25
$ws.Range("L132").Value = 40161
$ws.Range("M132").Value = 841.75
$ws.Range("N132").Value = -45221

$ws.Range("H136").Value = 1663123.4
$ws.Range("I136").Value = 1662054.5
$ws.Range("J136").Value = 1670784.1
$ws.Range("K136").Value = 4986163.5
$ws.Range("L136").Value = 5012352.300000001
$ws.Range("M136").Value = -4983613.5
$ws.Range("N136").Value = -5017452.300000001
